# The commit swaps the East-Asian default font from "DejaVu Sans" to
# "Tahoma" (docDefaults, Normal, Heading) and stamps an explicit
# complex-script font of "DejaVu Sans" onto a few styles that previously
# inherited it implicitly (List, Caption, Index).

$d = $word.ActiveDocument
$styles = $d.Styles

# docDefaults/rPrDefault covers the whole document as a fallback; Word's
# object model surfaces that through the "Normal" style together with the
# other named styles below.
$normal = $styles.Item("Normal")
$normal.Font.NameFarEast = "Tahoma"

$heading = $styles.Item("Heading")
$heading.Font.NameFarEast = "Tahoma"

$list = $styles.Item("List")
$list.Font.NameBi = "DejaVu Sans"

$caption = $styles.Item("Caption")
$caption.Font.NameBi = "DejaVu Sans"

$index = $styles.Item("Index")
$index.Font.NameBi = "DejaVu Sans"
